$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 23, shifting existing rows 23:42 down to 24:43.
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with data (same as the other rows in this
# block, but with its own date / price values).
$ws.Range("A23").Value = 5
$ws.Range("B23").Value = 'Macroferia Regional de Talca'
$ws.Range("C23").Value = 'Maule'
$ws.Range("D23").Value = 44762
$ws.Range("D23").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E23").Value = 7
$ws.Range("F23").Value = 100112043
$ws.Range("G23").Value = 'Pepino dulce'
$ws.Range("H23").Value = 'Cultivar IV Región'
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 16000
$ws.Range("L23").Value = 16000
$ws.Range("M23").Value = 16000
$ws.Range("N23").Value = '$/bandeja 18 kilos'
$ws.Range("O23").Value = 'Provincia de Limarí'
$ws.Range("P23").Value = 889
$ws.Range("Q23").Value = 18
$ws.Range("R23").Value = 'Hortaliza'
